$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C first: validation regex values (so shared strings land at
# indices 3-7, matching the target string table order)
$ws.Range("C2").Value = "[A-Z a-z].*"
$ws.Range("C3").Value = "[A-Z]{2}"
$ws.Range("C4").Value = "[0-9]{5}"
$ws.Range("C5").Value = "[0-9]{3}"
$ws.Range("C6").Value = "[A-Z]{1}"

# Column B: Xpath entries for the REST-style response fields (land at
# shared-string indices 8-12)
$ws.Range("B2").Value = "/Envelope/Body/GetInfoByStateResponse/GetInfoByStateResult/NewDataSet/Table[1]/CITY"
$ws.Range("B3").Value = "/Envelope/Body/GetInfoByStateResponse/GetInfoByStateResult/NewDataSet/Table[1]/STATE"
$ws.Range("B4").Value = "/Envelope/Body/GetInfoByStateResponse/GetInfoByStateResult/NewDataSet/Table[1]/ZIP"
$ws.Range("B5").Value = "/Envelope/Body/GetInfoByStateResponse/GetInfoByStateResult/NewDataSet/Table[1]/AREA_CODE"
$ws.Range("B6").Value = "/Envelope/Body/GetInfoByStateResponse/GetInfoByStateResult/NewDataSet/Table[1]/TIME_ZONE"

# Column widths (column B widened to fit the longer Xpath strings, new
# column C sized for the short regex values)
$ws.Columns.Item(2).ColumnWidth = 91.7
$ws.Columns.Item(3).ColumnWidth = 9.1666667

# Move selection to F4, matching the target worksheet view
$ws.Range("F4").Select()
